# Insert a new weekly price-report row for "Locoto" (row 118) and shift
# all subsequent rows (118-142) down by one, ending at row 143.
# This mirrors a new week of data being prepended into the existing
# chronological block of rows, pushing the rest of the rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 118; Excel shifts rows 118:142 down
# to 119:143 and copies formatting from the row above for the new row.
$ws.Rows.Item(118).Insert()

# Populate the new row 118 with the new week's data. Columns A, B, C, E,
# F, G, H, N, O, Q, R are constant across the whole data block, so copy
# them straight from the row below (119), which still holds the old
# (now shifted) row-118 data.
$ws.Range("A118").Value = $ws.Range("A119").Value()
$ws.Range("B118").Value = $ws.Range("B119").Value()
$ws.Range("C118").Value = $ws.Range("C119").Value()
$ws.Range("E118").Value = $ws.Range("E119").Value()
$ws.Range("F118").Value = $ws.Range("F119").Value()
$ws.Range("G118").Value = $ws.Range("G119").Value()
$ws.Range("H118").Value = $ws.Range("H119").Value()
$ws.Range("N118").Value = $ws.Range("N119").Value()
$ws.Range("O118").Value = $ws.Range("O119").Value()
$ws.Range("Q118").Value = $ws.Range("Q119").Value()
$ws.Range("R118").Value = $ws.Range("R119").Value()

# New-week specific values.
$ws.Range("D118").Value = 44855
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 160
$ws.Range("K118").Value = 19000
$ws.Range("L118").Value = 20000
$ws.Range("M118").Value = 19500
$ws.Range("P118").Value = 975
